$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    ,@(83, 2, "t-bone (angled0", "Possibly", "Yes", "No", "Yes", "sedan, suv", "No", "Night", "No", "collision happened at 2 lane road")
    ,@(84, 2, "rear-end", "No", "No", "No", "Yes", "bus, van", "No", "Day", "No", "happened at congested area, so there could be more involved in the accident")
    ,@(85, 1, "cause unknown", "No", "No", "No", "No", "scooter?", "No", "Day", "Yes", "Although the cause is unknown, it seems that a bike toppled over at an intersection")
    ,@(86, 2, "T-bone", "Yes", "Yes", "No", "Yes", "bus, sedan", "No", "Night", "No", "There don't seem to be many streetlights here, so it was hard to see specific vehicles")
    ,@(87, 2, "rear-end", "No", "No", "No", "Yes", "sedan, scooter", "No", "Day", "No", "Nobody seemed to be hurt by the end. Happened close to an intersection")
    ,@(88, 2, "T-bone", "Yes", "Yes", "No", "Yes", "suv x2", "No", "Day", "No", "Happened at a busy intersection, mild impact")
    ,@(89, 1, "topple", "No", "No", "No", "No", "motorbike", "No", "Day", "No", "Simple topple at an intersection, nothing to be done")
    ,@(90, 2, "light t-bone", "Yes", "Yes", "No", "Yes", "bus, motorbike", "No", "Day", "No", "Although light impact, motorcyclist probably needs medical attention")
    ,@(91, 2, "light t-bone", "No", "No", "No", "Yes", "buggy?, motorbike", "No", "Day", "No", "It doesn't look like anyone was hurt, but send over police just in case")
    ,@(92, 1, "topple", "Yes", "Yes", "No", "Yes", "18-wheeler", "No", "Day", "No", "18-wheeler took a sharp turn and toppled - debris present and damage to city inflicted")
    ,@(93, 1, "rear-end", "No", "No", "No", "No", "bike, truck", "No", "Day", "No", "Bike hits truck lightly. Cyclist can like just get up and move on")
    ,@(94, 0, "N/A", "No", "No", "No", "No", "N/A", "No", "Day, clear", "No", "No accident occurred")
    ,@(95, "2 (likely more)", "topple", "Yes", "Yes", "Yes", "Yes", "18-wheeler, sedan", "No", "Day, cloudy", "No", "18-wheeler topples on interstate and hits other vehicles in the process. Debris present. Maximum assistance is needed")
    ,@(96, "2 (likely more)", "topple", "Yes", "Yes", "Yes", "Yes", "18-wheeler, sedan", "No", "Day, cloudy", "No", "18-wheeler (again) topples on interstate and hits other vehicles in the process. Debris present. Maximum assistance is needed")
    ,@(97, 2, "light t-bone", "No", "No", "No", "Yes", "18-wheeler, van", "No", "Day, rainy", "No", "happened at intersection, light impact")
    ,@(98, 2, "light t-bone", "Possibly", "Yes", "No", "Yes", "bike, suv", "No", "Rainy, cloudy day", "No", "Suv hits cyclist - cyclist could be injured so medical assistance is necessary")
    ,@(99, 2, "t-bone (angled)", "No", "No", "No", "Yes", "suv, bus", "No", "Night", "No", "suv ran a red light and hit the bus lightly")
    ,@(100, 2, "rear-end", "No", "No", "No", "Yes", "suv x2", "No", "Rainy, cloudy day", "No", "light/mild rear end collision near a traffic light")
    ,@(101, 2, "t-bone", "Yes", "Yes", "No", "Yes", "bike, 18-wheeler?", "No", "Cloudy day?", "No", "big vehicle looks to have run a red light and, as a result, the biker hits the vehicle")
)

$startRow = 85
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

try {
    $excel.ActiveWindow.ScrollRow = 73
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
}

$ws.Range("L103").Select()
